try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item(1)
    $ws.Unprotect()
    $ws.Range("Z100").Value = "helloZ100"
    $ws.Protect($null, $true, $true, $true)
    Write-Output "set ok"
} catch {
    Write-Output "ERR: $_"
}
